# Add a new credential row (himanshu.sharma@diaspark.com / test123) below
# the existing rahul.sharma@nytimes.com row, matching the formatting that
# Excel applies to the existing data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "himanshu.sharma@diaspark.com"
$ws.Range("B3").Value = "test123"

# Re-apply the theme font color across both the existing and new data rows so
# they share a single (new) cell style, just like the source workbook does.
$ws.Range("A2:B3").Font.ThemeColor = 1
